$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntn1"
$ws.Cells.Item(2,3).Value = "Unc5b"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.990837
$ws.Cells.Item(2,8).Value = 5.972511000000001
$ws.Cells.Item(2,9).Value = 0.1122845585713437
$ws.Cells.Item(2,10).Value = 0.1122845585713437
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 9.581373333333334
$ws.Cells.Item(2,14).Value = 28.74412
$ws.Cells.Item(2,15).Value = 0.6094546925631529
$ws.Cells.Item(2,16).Value = 0.6094546925631529
$ws.Cells.Item(2,17).Value = 19.07495254281334
$ws.Cells.Item(2,18).Value = 171.67457288532
$ws.Cells.Item(2,19).Value = 0.06843235112368763
$ws.Cells.Item(2,20).Value = 0.06843235112368763

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntn1"
$ws.Cells.Item(3,3).Value = "Unc5b"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.990837
$ws.Cells.Item(3,8).Value = 5.972511000000001
$ws.Cells.Item(3,9).Value = 0.1122845585713437
$ws.Cells.Item(3,10).Value = 0.1122845585713437
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.468510333333334
$ws.Cells.Item(3,14).Value = 13.405531
$ws.Cells.Item(3,15).Value = 0.284234263364153
$ws.Cells.Item(3,16).Value = 0.2842342633641529
$ws.Cells.Item(3,17).Value = 8.896075706482335
$ws.Cells.Item(3,18).Value = 80.064681358341
$ws.Cells.Item(3,19).Value = 0.03191511879269498
$ws.Cells.Item(3,20).Value = 0.03191511879269497

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntn1"
$ws.Cells.Item(4,3).Value = "Unc5b"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.990837
$ws.Cells.Item(4,8).Value = 5.972511000000001
$ws.Cells.Item(4,9).Value = 0.1122845585713437
$ws.Cells.Item(4,10).Value = 0.1122845585713437
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.15426
$ws.Cells.Item(4,14).Value = 0.46278
$ws.Cells.Item(4,15).Value = 0.00981221351095027
$ws.Cells.Item(4,16).Value = 0.009812213510950268
$ws.Cells.Item(4,17).Value = 0.30710651562
$ws.Cells.Item(4,18).Value = 2.76395864058
$ws.Cells.Item(4,19).Value = 0.001101760062684826
$ws.Cells.Item(4,20).Value = 0.001101760062684826

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Ntn1"
$ws.Cells.Item(5,3).Value = "Unc5b"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.990837
$ws.Cells.Item(5,8).Value = 5.972511000000001
$ws.Cells.Item(5,9).Value = 0.1122845585713437
$ws.Cells.Item(5,10).Value = 0.1122845585713437
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.517079666666667
$ws.Cells.Item(5,14).Value = 4.551239
$ws.Cells.Item(5,15).Value = 0.09649883056174381
$ws.Cells.Item(5,16).Value = 0.0964988305617438
$ws.Cells.Item(5,17).Value = 3.020258332347667
$ws.Cells.Item(5,18).Value = 27.182324991129
$ws.Cells.Item(5,19).Value = 0.0108353285922763
$ws.Cells.Item(5,20).Value = 0.0108353285922763

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ntn1"
$ws.Cells.Item(6,3).Value = "Unc5b"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 11.42765333333333
$ws.Cells.Item(6,8).Value = 34.28296
$ws.Cells.Item(6,9).Value = 0.6445274073365515
$ws.Cells.Item(6,10).Value = 0.6445274073365515
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 9.581373333333334
$ws.Cells.Item(6,14).Value = 28.74412
$ws.Cells.Item(6,15).Value = 0.6094546925631529
$ws.Cells.Item(6,16).Value = 0.6094546925631529
$ws.Cells.Item(6,17).Value = 109.4926129105778
$ws.Cells.Item(6,18).Value = 985.4335161952001
$ws.Cells.Item(6,19).Value = 0.392810252886824
$ws.Cells.Item(6,20).Value = 0.392810252886824

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntn1"
$ws.Cells.Item(7,3).Value = "Unc5b"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 11.42765333333333
$ws.Cells.Item(7,8).Value = 34.28296
$ws.Cells.Item(7,9).Value = 0.6445274073365515
$ws.Cells.Item(7,10).Value = 0.6445274073365515
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.468510333333334
$ws.Cells.Item(7,14).Value = 13.405531
$ws.Cells.Item(7,15).Value = 0.284234263364153
$ws.Cells.Item(7,16).Value = 0.2842342633641529
$ws.Cells.Item(7,17).Value = 51.06458700575111
$ws.Cells.Item(7,18).Value = 459.58128305176
$ws.Cells.Item(7,19).Value = 0.1831967728423121
$ws.Cells.Item(7,20).Value = 0.183196772842312

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Ntn1"
$ws.Cells.Item(8,3).Value = "Unc5b"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.42765333333333
$ws.Cells.Item(8,8).Value = 34.28296
$ws.Cells.Item(8,9).Value = 0.6445274073365515
$ws.Cells.Item(8,10).Value = 0.6445274073365515
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.15426
$ws.Cells.Item(8,14).Value = 0.46278
$ws.Cells.Item(8,15).Value = 0.00981221351095027
$ws.Cells.Item(8,16).Value = 0.009812213510950268
$ws.Cells.Item(8,17).Value = 1.7628298032
$ws.Cells.Item(8,18).Value = 15.8654682288
$ws.Cells.Item(8,19).Value = 0.006324240534445458
$ws.Cells.Item(8,20).Value = 0.006324240534445458

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Ntn1"
$ws.Cells.Item(9,3).Value = "Unc5b"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.42765333333333
$ws.Cells.Item(9,8).Value = 34.28296
$ws.Cells.Item(9,9).Value = 0.6445274073365515
$ws.Cells.Item(9,10).Value = 0.6445274073365515
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.517079666666667
$ws.Cells.Item(9,14).Value = 4.551239
$ws.Cells.Item(9,15).Value = 0.09649883056174381
$ws.Cells.Item(9,16).Value = 0.0964988305617438
$ws.Cells.Item(9,17).Value = 17.33666050971556
$ws.Cells.Item(9,18).Value = 156.02994458744
$ws.Cells.Item(9,19).Value = 0.06219614107296991
$ws.Cells.Item(9,20).Value = 0.0621961410729699

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ntn1"
$ws.Cells.Item(10,3).Value = "Unc5b"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.2266433333333333
$ws.Cells.Item(10,8).Value = 0.67993
$ws.Cells.Item(10,9).Value = 0.01278283789002879
$ws.Cells.Item(10,10).Value = 0.01278283789002879
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 9.581373333333334
$ws.Cells.Item(10,14).Value = 28.74412
$ws.Cells.Item(10,15).Value = 0.6094546925631529
$ws.Cells.Item(10,16).Value = 0.6094546925631529
$ws.Cells.Item(10,17).Value = 2.171554390177778
$ws.Cells.Item(10,18).Value = 19.5439895116
$ws.Cells.Item(10,19).Value = 0.007790560536352118
$ws.Cells.Item(10,20).Value = 0.007790560536352118

$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Ntn1"
$ws.Cells.Item(11,3).Value = "Unc5b"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.2266433333333333
$ws.Cells.Item(11,8).Value = 0.67993
$ws.Cells.Item(11,9).Value = 0.01278283789002879
$ws.Cells.Item(11,10).Value = 0.01278283789002879
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 4.468510333333334
$ws.Cells.Item(11,14).Value = 13.405531
$ws.Cells.Item(11,15).Value = 0.284234263364153
$ws.Cells.Item(11,16).Value = 0.2842342633641529
$ws.Cells.Item(11,17).Value = 1.012758076981111
$ws.Cells.Item(11,18).Value = 9.11482269283
$ws.Cells.Item(11,19).Value = 0.003633320511375717
$ws.Cells.Item(11,20).Value = 0.003633320511375716

$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Ntn1"
$ws.Cells.Item(12,3).Value = "Unc5b"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.2266433333333333
$ws.Cells.Item(12,8).Value = 0.67993
$ws.Cells.Item(12,9).Value = 0.01278283789002879
$ws.Cells.Item(12,10).Value = 0.01278283789002879
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.15426
$ws.Cells.Item(12,14).Value = 0.46278
$ws.Cells.Item(12,15).Value = 0.00981221351095027
$ws.Cells.Item(12,16).Value = 0.009812213510950268
$ws.Cells.Item(12,17).Value = 0.03496200059999999
$ws.Cells.Item(12,18).Value = 0.3146580054
$ws.Cells.Item(12,19).Value = 0.0001254279346528275
$ws.Cells.Item(12,20).Value = 0.0001254279346528275

$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Ntn1"
$ws.Cells.Item(13,3).Value = "Unc5b"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.2266433333333333
$ws.Cells.Item(13,8).Value = 0.67993
$ws.Cells.Item(13,9).Value = 0.01278283789002879
$ws.Cells.Item(13,10).Value = 0.01278283789002879
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.517079666666667
$ws.Cells.Item(13,14).Value = 4.551239
$ws.Cells.Item(13,15).Value = 0.09649883056174381
$ws.Cells.Item(13,16).Value = 0.0964988305617438
$ws.Cells.Item(13,17).Value = 0.3438359925855555
$ws.Cells.Item(13,18).Value = 3.09452393327
$ws.Cells.Item(13,19).Value = 0.001233528907648127
$ws.Cells.Item(13,20).Value = 0.001233528907648127

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Ntn1"
$ws.Cells.Item(14,3).Value = "Unc5b"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.085149333333333
$ws.Cells.Item(14,8).Value = 12.255448
$ws.Cells.Item(14,9).Value = 0.230405196202076
$ws.Cells.Item(14,10).Value = 0.230405196202076
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 9.581373333333334
$ws.Cells.Item(14,14).Value = 28.74412
$ws.Cells.Item(14,15).Value = 0.6094546925631529
$ws.Cells.Item(14,16).Value = 0.6094546925631529
$ws.Cells.Item(14,17).Value = 39.14134088508444
$ws.Cells.Item(14,18).Value = 352.27206796576
$ws.Cells.Item(14,19).Value = 0.1404215280162892
$ws.Cells.Item(14,20).Value = 0.1404215280162892

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Ntn1"
$ws.Cells.Item(15,3).Value = "Unc5b"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.085149333333333
$ws.Cells.Item(15,8).Value = 12.255448
$ws.Cells.Item(15,9).Value = 0.230405196202076
$ws.Cells.Item(15,10).Value = 0.230405196202076
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 4.468510333333334
$ws.Cells.Item(15,14).Value = 13.405531
$ws.Cells.Item(15,15).Value = 0.284234263364153
$ws.Cells.Item(15,16).Value = 0.2842342633641529
$ws.Cells.Item(15,17).Value = 18.25453200920978
$ws.Cells.Item(15,18).Value = 164.290788082888
$ws.Cells.Item(15,19).Value = 0.06548905121777021
$ws.Cells.Item(15,20).Value = 0.0654890512177702

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Ntn1"
$ws.Cells.Item(16,3).Value = "Unc5b"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.085149333333333
$ws.Cells.Item(16,8).Value = 12.255448
$ws.Cells.Item(16,9).Value = 0.230405196202076
$ws.Cells.Item(16,10).Value = 0.230405196202076
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.15426
$ws.Cells.Item(16,14).Value = 0.46278
$ws.Cells.Item(16,15).Value = 0.00981221351095027
$ws.Cells.Item(16,16).Value = 0.009812213510950268
$ws.Cells.Item(16,17).Value = 0.6301751361599999
$ws.Cells.Item(16,18).Value = 5.671576225439999
$ws.Cells.Item(16,19).Value = 0.002260784979167158
$ws.Cells.Item(16,20).Value = 0.002260784979167158

$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Ntn1"
$ws.Cells.Item(17,3).Value = "Unc5b"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.085149333333333
$ws.Cells.Item(17,8).Value = 12.255448
$ws.Cells.Item(17,9).Value = 0.230405196202076
$ws.Cells.Item(17,10).Value = 0.230405196202076
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 1.517079666666667
$ws.Cells.Item(17,14).Value = 4.551239
$ws.Cells.Item(17,15).Value = 0.09649883056174381
$ws.Cells.Item(17,16).Value = 0.0964988305617438
$ws.Cells.Item(17,17).Value = 6.197496988896888
$ws.Cells.Item(17,18).Value = 55.777472900072
$ws.Cells.Item(17,19).Value = 0.02223383198884947
$ws.Cells.Item(17,20).Value = 0.02223383198884947

